$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The original blank placeholder row 9 (C9/E9/F9/G9, all empty) was removed from the
# expense table. Deleting the whole row shifts every row below it up by one and lets
# Excel auto-adjust the trailing SUM() formula range, which matches the target sheet
# exactly (old row 42 "合计" becomes row 41, SUM(H2:H41) becomes SUM(H2:H40)).
$ws.Rows("9:9").Delete()

# Row 7: new purchase entry - Raspberry Pi kit bought via Taobao.
$ws.Range("B7").Value = 20180820
$ws.Range("C7").Value = "元器件"
$ws.Range("G11").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D11").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D7").Value = "树莓派1套"
$ws.Range("E7").Value = "淘宝"
$ws.Range("F7").Value = "https://item.taobao.com/item.htm?spm=a1z09.2.0.0.2b962e8dlwdNUv&id=527576110046&_u=e3s1sni8961"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 280

# Row 8: parking-fee entry at Liandong U Valley.
$ws.Range("B8").Value = 20180913
$ws.Range("D8").Value = "联东U谷停车费"
$ws.Range("H8").Value = 17.5

# Restore the view to the top of the sheet, matching the saved selection state.
$ws.Range("F18").Select()
